# Apply "contingencies with rene fine": insert line7/line8 rows and shift
# the extr1..extr8 rows down by two, updating their contingency values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target state for rows 8..17 (A=index, B=name, C=from_bus, D=to_bus, E=in_service)
$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($r in $rows) {
    $row = $r.Row
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $r.A
    # Match the "index" column formatting used elsewhere in the sheet:
    # bold font, thin border all around, centered horizontally / top vertically.
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
}
